# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.403.06"
Set-TextValue "E2" "  +2.35%  "
Set-TextValue "D3" "2.227.42"
Set-TextValue "E3" "  +0.45%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "231.11"
Set-TextValue "E5" "  +0.58%  "
Set-TextValue "E6" "  -0.44%  "
Set-TextValue "D7" "60.66"
Set-TextValue "E7" "  -0.61%  "
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D10" "0.0903"
Set-TextValue "E10" "  +1.15%  "
Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  +0.27%  "
Set-TextValue "D12" "2.557.91"
Set-TextValue "E12" "  +0.50%  "
Set-TextValue "D13" "15.57"
Set-TextValue "E13" "  +0.10%  "
Set-TextValue "D14" "22.28"
Set-TextValue "E14" "  +3.64%  "
Set-TextValue "D15" "5.61"
Set-TextValue "E15" "  +1.52%  "
Set-TextValue "D16" "0.799"
Set-TextValue "E16" "  +0.69%  "
Set-TextValue "D17" "2.250.42"
Set-TextValue "E17" "  +1.62%  "
Set-TextValue "D18" "42.309.71"
Set-TextValue "E18" "  +2.45%  "
Set-TextValue "D19" "0.0₃0942"
Set-TextValue "E19" "  +5.44%  "
Set-TextValue "D20" "6.17"
Set-TextValue "E20" "  +2.39%  "
Set-TextValue "D21" "72.22"
Set-TextValue "E21" "  -0.50%  "
Set-TextValue "D22" "244.15"
Set-TextValue "E22" "  -2.89%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  -0.29%  "
Set-TextValue "D24" "2.38"
Set-TextValue "E24" "  +3.13%  "
Set-TextValue "D25" "2.39"
Set-TextValue "E25" "  +0.30%  "
Set-TextValue "D26" "9.72"
Set-TextValue "E26" "  +2.01%  "
Set-TextValue "D27" "169.52"
Set-TextValue "E27" "  +1.38%  "
Set-TextValue "D28" "0.141"
Set-TextValue "E28" "  +0.28%  "
Set-TextValue "D29" "20.33"
Set-TextValue "E29" "  +2.42%  "
Set-TextValue "E30" "  +2.32%  "
Set-TextValue "D31" "2.66"
Set-TextValue "E31" "  +4.85%  "
Set-TextValue "D32" "0.120"
Set-TextValue "E32" "  -1.42%  "
Set-TextValue "D33" "5.02"
Set-TextValue "E33" "  +0.82%  "
Set-TextValue "D34" "4.64"
Set-TextValue "D35" "0.0651"
Set-TextValue "E35" "  +5.26%  "
Set-TextValue "D36" "6.39"
Set-TextValue "E36" "  -1.97%  "
Set-TextValue "E37" "  -0.45%  "
Set-TextValue "E38" "  -3.55%  "
Set-TextValue "E39" "  +6.13%  "
Set-TextValue "D40" "0.999"
Set-TextValue "E40" "  -0.03%  "
Set-TextValue "D41" "0.000232"
Set-TextValue "E41" "  -1.53%  "
Set-TextValue "D42" "8.58"
Set-TextValue "E42" "  -0.33%  "
Set-TextValue "B43" "TrustWalletToken"
Set-TextValue "C43" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "1.21"
Set-TextValue "E43" "  +1.19%  "
Set-TextValue "B44" "Cronos"
Set-TextValue "C44" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0960"
Set-TextValue "E44" "  -1.71%  "
Set-TextValue "B45" "FTXToken"
Set-TextValue "C45" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D45" "4.42"
Set-TextValue "E45" "  -7.49%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "96.95"
Set-TextValue "E46" "  -1.90%  "
Set-TextValue "D47" "1.458.04"
Set-TextValue "E47" "  -0.37%  "
Set-TextValue "E48" "  -1.16%  "
Set-TextValue "D49" "1.08"
Set-TextValue "E49" "  +0.16%  "
Set-TextValue "D50" "16.00"
Set-TextValue "E50" "  -3.10%  "
Set-TextValue "D51" "2.23"
Set-TextValue "E51" "  +4.70%  "
